# Journal de travail - add new work-log entries (July 2023) and extend totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# --- Clear the old totals-row content; it moves down to rows 88/90 -------
$ws.Range("C77").ClearContents() | Out-Null
$ws.Range("C79").ClearContents() | Out-Null
$ws.Range("D79").ClearContents() | Out-Null

# --- New data rows (72-74) ---------------------------------------------
$ws.Range("A72").Value = 45118
$ws.Range("B72").Value = "Rédaction"
$ws.Range("C72").Value = 3
$ws.Range("D72").Value = "Rapport: plan et idées"

$ws.Range("A73").Value = 45118
$ws.Range("B73").Value = "Implémentation"
$ws.Range("C73").Value = 3
$ws.Range("D73").Value = "Debug"

$ws.Range("A74").Value = 45119
$ws.Range("B74").Value = "Implémentation"
$ws.Range("C74").Value = 2
$ws.Range("D74").Value = "Validation JSON schéma, "

# --- Extend the dated-row formatting down to row 86 ---------------------
# Column A keeps reusing the existing date-cell style (same as A2:A74) by
# copying the format from an already-styled cell - this preserves the
# workbook's original style index instead of minting a new one.
$ws.Range("A72").Copy() | Out-Null
$ws.Range("A75:A86").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column C (Temps [h]) for the blank rows gets an explicit "General"
# number format, matching the new style added to the workbook.
$ws.Range("C75:C85").NumberFormat = "General"

# --- Totals --------------------------------------------------------------
$ws.Range("C88").Formula = "=SUM(C2:C87)"

$ws.Range("C90").Value = "Temps plein:"
$ws.Range("D90").Formula = "=SUM(C54:C86)"

# --- Resize the structured table to include the new rows -----------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E88")) | Out-Null

# --- View state: keep selection in sync with the edited area -------------
$ws.Activate() | Out-Null
$ws.Range("D74").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
